# Monthly database update: drop the oldest month column (ماه 7 منتهی به 1397/07)
# and append a new month column (ماه 9 منتهی به 1401/09) at the end, shifting every
# month-indexed row one column to the left. Also bump the copyright year.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copyright year bump
$ws.Range("B3").Value = "Copyright @2015 - 2023"

# 2) Rows that carry one value per month in columns E(5)..BB(54).
#    For each of these rows we shift F..BB (cols 6..54) into E..BA (cols 5..53)
#    and then write a brand-new value into BB (col 54).
$firstCol = 5   # E
$lastCol  = 54  # BB

$newMonthLabel = "ماه 9 منتهی به 1401/09"

$rowsAndNewValues = @{
    8  = $newMonthLabel
    11 = 56886
    13 = 0
    14 = 56886
    18 = $newMonthLabel
    21 = 59420
    23 = 0
    25 = 0
    26 = 59420
    30 = $newMonthLabel
    33 = 5017871
    35 = 0
    37 = 0
    39 = 0
    40 = 5017871
    44 = $newMonthLabel
    47 = 84447509
}

foreach ($row in $rowsAndNewValues.Keys) {
    for ($col = $firstCol; $col -lt $lastCol; $col++) {
        $srcVal = $ws.Cells.Item($row, $col + 1).Value2
        $ws.Cells.Item($row, $col).Value = $srcVal
    }
    $ws.Cells.Item($row, $lastCol).Value = $rowsAndNewValues[$row]
}
